$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 980
$ws.Range("Y2").Value = 980
$ws.Range("AB2").Value = 980
$ws.Range("AC2").Value = 980
$ws.Range("AD2").Value = 980
$ws.Range("AG2").Value = 980
$ws.Range("AH2").Value = 980
$ws.Range("X3").Value = 980
$ws.Range("Y3").Value = 980
$ws.Range("AB3").Value = 980
$ws.Range("AC3").Value = 980
$ws.Range("AD3").Value = 980
$ws.Range("AG3").Value = 980
$ws.Range("AH3").Value = 980
$ws.Range("F4").Value = 1.49
$ws.Range("I4").Value = 19.5
$ws.Range("K4").Value = 7.8
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 1.76
$ws.Range("V4").Value = 1.05
$ws.Range("X4").Value = 980
$ws.Range("Y4").Value = 980
$ws.Range("AB4").Value = 980
$ws.Range("AC4").Value = 980
$ws.Range("AD4").Value = 980
$ws.Range("AG4").Value = 980
$ws.Range("AH4").Value = 980
$ws.Range("N5").Value = 1.02
$ws.Range("P5").Value = 1.02
$ws.Range("R5").Value = 1.02
$ws.Range("X5").Value = 980
$ws.Range("Y5").Value = 980
$ws.Range("AB5").Value = 980
$ws.Range("AC5").Value = 980
$ws.Range("AD5").Value = 980
$ws.Range("AG5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("L6").Value = 1.46
$ws.Range("T6").Value = 1.86
$ws.Range("H7").Value = 38
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 16.5
$ws.Range("K7").Value = 18
$ws.Range("S7").Value = 1.77
$ws.Range("T7").Value = 2.62
$ws.Range("U7").Value = 1.54
$ws.Range("W7").Value = 11.5
$ws.Range("Y7").Value = 980
$ws.Range("AC7").Value = 42
$ws.Range("AD7").Value = 980
$ws.Range("I8").Value = 7.6
$ws.Range("J8").Value = 4.6
$ws.Range("H9").Value = 1.7
$ws.Range("X9").Value = 980
$ws.Range("Y9").Value = 980
$ws.Range("AB9").Value = 980
$ws.Range("AC9").Value = 980
$ws.Range("AD9").Value = 980
$ws.Range("AG9").Value = 980
$ws.Range("AH9").Value = 980
$ws.Range("R11").Value = 1.2
$ws.Range("L12").Value = 1.49
$ws.Range("M12").Value = 1.09
$ws.Range("H13").Value = 3.3
$ws.Range("I13").Value = 3.35
$ws.Range("L13").Value = 1.39
$ws.Range("Q13").Value = 1.89
$ws.Range("X13").Value = 14.5
$ws.Range("AC13").Value = 8
$ws.Range("G14").Value = 2.42
$ws.Range("L14").Value = 1.38
$ws.Range("W14").Value = 1.7
$ws.Range("AB14").Value = 9.4
$ws.Range("H15").Value = 9
$ws.Range("I15").Value = 9.4
$ws.Range("L15").Value = 1.25
$ws.Range("N15").Value = 5.1
$ws.Range("V15").Value = 1.11
$ws.Range("AA15").Value = 290
$ws.Range("AG15").Value = 9.800000000000001
$ws.Range("AH15").Value = 24
$ws.Range("AO15").Value = 120
$ws.Range("T16").Value = 1.66
$ws.Range("X16").Value = 60
$ws.Range("AN16").Value = 3.7
$ws.Range("P17").Value = 1.87
$ws.Range("T17").Value = 1.88
$ws.Range("H18").Value = 2.08
$ws.Range("I18").Value = 2.1
$ws.Range("Q18").Value = 1.45
$ws.Range("V18").Value = 1.9
$ws.Range("X18").Value = 32
$ws.Range("Y18").Value = 18.5
$ws.Range("I19").Value = 6.4
$ws.Range("V19").Value = 1.18
$ws.Range("Y19").Value = 22
$ws.Range("AF19").Value = 9.4
$ws.Range("H20").Value = 13
$ws.Range("J20").Value = 6.4
$ws.Range("Q20").Value = 1.52
$ws.Range("AN20").Value = 4.1
$ws.Range("F21").Value = 1.35
$ws.Range("G21").Value = 1.36
$ws.Range("H21").Value = 11
$ws.Range("I21").Value = 11.5
$ws.Range("J21").Value = 5.8
$ws.Range("K21").Value = 6
$ws.Range("P21").Value = 2.4
$ws.Range("R21").Value = 1.55
$ws.Range("S21").Value = 2.68
$ws.Range("W21").Value = 3.75
$ws.Range("AA21").Value = 420
$ws.Range("AD21").Value = 40
$ws.Range("AE21").Value = 180
$ws.Range("AF21").Value = 8.199999999999999
$ws.Range("AI21").Value = 150
$ws.Range("AJ21").Value = 10.5
$ws.Range("AO21").Value = 210
$ws.Range("X22").Value = 980
$ws.Range("Y22").Value = 980
$ws.Range("AB22").Value = 980
$ws.Range("AC22").Value = 980
$ws.Range("AD22").Value = 980
$ws.Range("AG22").Value = 980
$ws.Range("AH22").Value = 980
$ws.Range("H23").Value = 2.72
$ws.Range("R23").Value = 1.23
$ws.Range("S23").Value = 2.2
$ws.Range("X23").Value = 980
$ws.Range("Y23").Value = 980
$ws.Range("AB23").Value = 980
$ws.Range("AC23").Value = 980
$ws.Range("AD23").Value = 980
$ws.Range("AG23").Value = 980
$ws.Range("AH23").Value = 980
